# Rename the inline logo pictures that live in the document's headers and
# footers:
#   - the two "BTec_Logo-Orange" pictures (in the first-page header and the
#     default header) go from "image1.jpg" to "image2.jpg"
#   - the two Pearson Edexcel logo pictures (in the first-page footer and
#     the default footer) go from "image2.png" to "image1.png"
#
# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2 — the Headers /
# Footers collections are indexed by those constants, not sequentially.

$d = $word.ActiveDocument

function Rename-LogoShape($range, $newName) {
    if ($range.InlineShapes.Count -ge 1) {
        $shape = $range.InlineShapes(1)
        $shape.Name = $newName
    }
}

foreach ($sec in $d.Sections) {
    # Headers: BTec_Logo-Orange image1.jpg -> image2.jpg
    $hdrPrimary = $sec.Headers(1)
    if ($hdrPrimary.Exists) { Rename-LogoShape $hdrPrimary.Range "image2.jpg" }

    $hdrFirst = $sec.Headers(2)
    if ($hdrFirst.Exists) { Rename-LogoShape $hdrFirst.Range "image2.jpg" }

    # Footers: Pearson logo image2.png -> image1.png
    $ftrPrimary = $sec.Footers(1)
    if ($ftrPrimary.Exists) { Rename-LogoShape $ftrPrimary.Range "image1.png" }

    $ftrFirst = $sec.Footers(2)
    if ($ftrFirst.Exists) { Rename-LogoShape $ftrFirst.Range "image1.png" }
}
